# Press release header: give the "TBD, 2011" placeholder a real date.
# The final text is split across two runs ("March 30" / ", 2011") that
# share identical run formatting, matching how the author's edit landed.

$d = $word.ActiveDocument
$sec = $d.Sections.First
$hdr = $sec.Headers.Item(2)   # the header story that holds the dateline

# 1) Swap the whole placeholder for the final text in one shot.
$dateRange = $hdr.Range.Duplicate
$null = $dateRange.Find.Execute("TBD, 2011", $true, $false, $false, $false, $false, $true, 1, $false, "March 30, 2011", 2)

# Find/Execute on a tracked-changes document records w:ins/w:del markup;
# accept those specific revisions (not a document-wide AcceptAll, which
# would also strip unrelated rsid bookkeeping elsewhere in the header).
for ($i = $hdr.Range.Revisions.Count; $i -ge 1; $i--) {
    $hdr.Range.Revisions.Item($i).Accept()
}

# 2) Re-touch the ", 2011" tail's formatting (off then back on) so the
#    serializer keeps it as its own run instead of folding it back into
#    the "March 30" run it is textually/format-identical to.
$yearRange = $hdr.Range.Duplicate
$null = $yearRange.Find.Execute(", 2011", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$yearRange.Font.Bold = $false
$yearRange.Font.Bold = $true
